# This script applies the textual and formatting edits described by the
# provided unified diff to the active Word document.

$d = $word.ActiveDocument

function Replace-Exact {
    param(
        [string]$OldText,
        [string]$NewText,
        [bool]$Bold
    )

    $rng = $d.Content
    $found = $rng.Find.Execute($OldText, $true, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text"
    }
    if ($Bold) {
        $rng.Font.Bold = 1
    }
    $rng.Text = $NewText
}

# 1. Title run: turn bold on and change the title text.
Replace-Exact "디자인 팀의 주요 업무 문서" "팀 핵심 책임 문서 디자인" $true

# 2. Intro paragraph: "그래픽 디자인 연구소" -> "Graphic Design Institute"
Replace-Exact ": 이 문서에서는 그래픽 디자인 연구소의 모든 디자인 팀 구성원의 핵심 책임을 간략하게 설명합니다." `
              ": 이 문서에서는 Graphic Design Institute의 모든 디자인 팀 구성원의 핵심 책임을 간략하게 설명합니다." `
              $false

# 3. Designer collaboration responsibility bullet.
Replace-Exact ": 다른 디자이너, 개발자 및 이해 관계자와 협력하여 프로젝트 요구 사항을 충족하는 고품질 디자인을 만듭니다. " `
              ": 다른 디자이너, 개발자, 이해 관계자들과 공동으로 작업을 진행하여 프로젝트 요구 사항을 충족하는 고품질 디자인을 만들어야 합니다. " `
              $false

# 4. Visual design responsibility bullet.
Replace-Exact ": 사용자에게 친숙하고 접근성이 뛰어나며 반응성이 뛰어난 시각적으로 매력적인 디자인을 만듭니다. " `
              ": 사용자들이 익숙한 방식으로 쉽게 활용할 수 있으며 적극적으로 반응할 수 있는 멋진 스타일의 디자인을 만들어야 합니다. " `
              $false

# 5. Communication responsibility bullet.
Replace-Exact ": 팀 구성원, 이해 관계자 및 클라이언트와 효과적으로 통신하여 프로젝트 요구 사항을 충족하는지 확인합니다. " `
              ": 프로젝트 요구 사항을 충족할 수 있도록 팀 구성원, 이해 관계자, 고객과 효율적으로 커뮤니케이션합니다. " `
              $false

# 6. "연구" (bold heading run) -> "리서치"
Replace-Exact "연구" "리서치" $false

# 7. Research responsibility sentence: "연구를 수행합니다" -> "리서치를 수행합니다"
Replace-Exact ": 디자인 결정을 알리기 위해 사용자 요구 사항, 기본 설정 및 동작을 식별하기 위한 연구를 수행합니다. " `
              ": 디자인 결정을 알리기 위해 사용자 요구 사항, 기본 설정 및 동작을 식별하기 위한 리서치를 수행합니다. " `
              $false

# 8. Usability testing responsibility bullet.
Replace-Exact ": 디자인이 사용자 요구를 충족하고 모든 사용자가 액세스할 수 있도록 유용성 테스트를 수행합니다. " `
              ": 사용 편의성 테스트를 수행하여 디자인이 사용자의 요구를 충족하며 모든 사용자가 쉽게 사용 가능한 상태인지를 확인해야 합니다. " `
              $false

# 9. "전문 개발" -> "전문적인 개발"
Replace-Exact "전문 개발" "전문적인 개발" $false
